$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge the BOUNDARY header cells first so uniform formatting applies cleanly after ---
$ws.Range("AS1:AZ1").Merge()

# --- Formatting: replicate header/epsilon style (bold, centered, top-aligned, thin border) ---
$hdr = $ws.Range("AS1:AZ2")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# --- Row 1: BOUNDARY header label (rest of merged cells stay blank) ---
$ws.Range("AS1").Value = "BOUNDARY"

# --- Row 2: epsilon values stored as text (quote-prefixed to avoid numeric coercion) ---
$row2 = New-Object 'object[,]' 1,8
$row2[0,0] = "'0.01"
$row2[0,1] = "'0.02"
$row2[0,2] = "'0.03"
$row2[0,3] = "'0.04"
$row2[0,4] = "'0.05"
$row2[0,5] = "'0.07"
$row2[0,6] = "'0.10"
$row2[0,7] = "'0.20"
$ws.Range("AS2:AZ2").Value = $row2

# --- Rows 4-12 (skip 3): numeric attack-result data ---
$r4 = New-Object 'object[,]' 1,8
$r4[0,0] = 4.186208724975586
$r4[0,1] = 4.249595165252686
$r4[0,2] = 4.302657604217529
$r4[0,3] = 4.353903293609619
$r4[0,4] = 4.388765335083008
$r4[0,5] = 4.730515480041504
$r4[0,6] = 5.006624221801758
$r4[0,7] = 7.04196310043335
$ws.Range("AS4:AZ4").Value = $r4

$r5 = New-Object 'object[,]' 1,8
$r5[0,0] = 5.322792225955141
$r5[0,1] = 5.364417903781465
$r5[0,2] = 5.440529341188144
$r5[0,3] = 5.47907246362425
$r5[0,4] = 5.538794076546599
$r5[0,5] = 5.940059391328488
$r5[0,6] = 6.324484750621179
$r5[0,7] = 8.787343581911573
$ws.Range("AS5:AZ5").Value = $r5

$r6 = New-Object 'object[,]' 1,8
$r6[0,0] = 0.9996174573898315
$r6[0,1] = 0.999610960483551
$r6[0,2] = 0.9995973110198975
$r6[0,3] = 0.9996079206466675
$r6[0,4] = 0.9995827674865723
$r6[0,5] = 0.9994970560073853
$r6[0,6] = 0.9994394183158875
$r6[0,7] = 0.9989516735076904
$ws.Range("AS6:AZ6").Value = $r6

$r7 = New-Object 'object[,]' 1,8
$r7[0,0] = 2.768602848052979
$r7[0,1] = 2.931120872497559
$r7[0,2] = 3.079442024230957
$r7[0,3] = 3.175944566726685
$r7[0,4] = 3.539762020111084
$r7[0,5] = 4.408213138580322
$r7[0,6] = 5.126472949981689
$r7[0,7] = 8.938190460205078
$ws.Range("AS7:AZ7").Value = $r7

$r8 = New-Object 'object[,]' 1,8
$r8[0,0] = 3.709892173907027
$r8[0,1] = 3.876880312603415
$r8[0,2] = 4.108211918417303
$r8[0,3] = 4.086089909995064
$r8[0,4] = 4.608258024753505
$r8[0,5] = 5.58454845737579
$r8[0,6] = 6.494436964573326
$r8[0,7] = 11.17392798280629
$ws.Range("AS8:AZ8").Value = $r8

$r9 = New-Object 'object[,]' 1,8
$r9[0,0] = 0.9997526407241821
$r9[0,1] = 0.9997299909591675
$r9[0,2] = 0.9996969699859619
$r9[0,3] = 0.9996999502182007
$r9[0,4] = 0.9996183514595032
$r9[0,5] = 0.9994405508041382
$r9[0,6] = 0.9992414116859436
$r9[0,7] = 0.9977531433105469
$ws.Range("AS9:AZ9").Value = $r9

$r10 = New-Object 'object[,]' 1,8
$r10[0,0] = 2.971150875091553
$r10[0,1] = 3.163308620452881
$r10[0,2] = 3.416063070297241
$r10[0,3] = 3.897083282470703
$r10[0,4] = 4.299318313598633
$r10[0,5] = 5.331329822540283
$r10[0,6] = 6.975700855255127
$r10[0,7] = 11.95866680145264
$ws.Range("AS10:AZ10").Value = $r10

$r11 = New-Object 'object[,]' 1,8
$r11[0,0] = 3.831296891241326
$r11[0,1] = 4.097779365777375
$r11[0,2] = 4.376790361557034
$r11[0,3] = 4.945704922164688
$r11[0,4] = 5.395146428084112
$r11[0,5] = 6.785276843433952
$r11[0,6] = 8.638624404614655
$r11[0,7] = 15.06382600004944
$ws.Range("AS11:AZ11").Value = $r11

$r12 = New-Object 'object[,]' 1,8
$r12[0,0] = 0.9998003244400024
$r12[0,1] = 0.9997552037239075
$r12[0,2] = 0.9997143149375916
$r12[0,3] = 0.9996235966682434
$r12[0,4] = 0.999544084072113
$r12[0,5] = 0.9992557764053345
$r12[0,6] = 0.9987260103225708
$r12[0,7] = 0.9959961771965027
$ws.Range("AS12:AZ12").Value = $r12
